$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text format first so that
# numeric-looking values ("332.69", "1.26%", etc.) are stored as
# literal text (matching the original inlineStr cells) instead of
# being auto-converted to numbers/percentages by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '332.69'
$ws.Range('E2').Value = '1.26%'
$ws.Range('D3').Value = '44.00'
$ws.Range('E3').Value = '5.89%'
$ws.Range('D4').Value = '5.870'
$ws.Range('E4').Value = '4.48%'
$ws.Range('D5').Value = '0.08341'
$ws.Range('E5').Value = '1.96%'
$ws.Range('D6').Value = '8.791'
$ws.Range('E6').Value = '0.60%'
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').Value = '1.969'
$ws.Range('E7').Value = '-2.39%'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').Value = '2.898'
$ws.Range('E8').Value = '-1.82%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '0.9359'
$ws.Range('E9').Value = '1.46%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '0.1260'
$ws.Range('E10').Value = '-1.27%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '0.1945'
$ws.Range('E11').Value = '-0.51%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.09465'
$ws.Range('E12').Value = '1.34%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '0.04267'
$ws.Range('E13').Value = '11.96%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '0.1067'
$ws.Range('E14').Value = '0.65%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '0.001303'
$ws.Range('E15').Value = '-0.16%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.005923'
$ws.Range('E16').Value = '-3.92%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.504'
$ws.Range('E17').Value = '1.74%'
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = '4.502'
$ws.Range('E18').Value = '-0.49%'
$ws.Range('E19').Value = '0.74%'
$ws.Range('D20').Value = '8.801'
$ws.Range('E20').Value = '6.09%'
$ws.Range('D21').Value = '0.1372'
$ws.Range('E21').Value = '-0.55%'
$ws.Range('D22').Value = '0.2633'
$ws.Range('E22').Value = '9.26%'
$ws.Range('D23').Value = '0.04423'
$ws.Range('E23').Value = '0.17%'
$ws.Range('D24').Value = '0.001258'
$ws.Range('E24').Value = '0.04%'
$ws.Range('D25').Value = '0.004410'
$ws.Range('E25').Value = '2.06%'
$ws.Range('D26').Value = '0.0001191'
$ws.Range('E26').Value = '0.76%'
$ws.Range('D27').Value = '0.0003994'
$ws.Range('D39').Value = '0.02816'
$ws.Range('E39').Value = '1.66%'
$ws.Range('D40').Value = '0.05697'
$ws.Range('E40').Value = '5.35%'
$ws.Range('D41').Value = '0.007907'
$ws.Range('E41').Value = '2.70%'
$ws.Range('D42').Value = '0.1428'
$ws.Range('E42').Value = '0.67%'
$ws.Range('D43').Value = '0.009043'
$ws.Range('E43').Value = '0.88%'
$ws.Range('D44').Value = '0.002157'
$ws.Range('E44').Value = '0.63%'
$ws.Range('E45').Value = '-9.81%'
$ws.Range('D46').Value = '0.00007194'
$ws.Range('E46').Value = '9.80%'
$ws.Range('D47').Value = '0.00000000751'
$ws.Range('E47').Value = '-0.03%'
$ws.Range('D48').Value = '0.003247'
$ws.Range('E48').Value = '1.36%'
$ws.Range('E49').Value = '-0.01%'
$ws.Range('D50').Value = '0.00002103'
$ws.Range('E50').Value = '-0.03%'
$ws.Range('D51').Value = '0.0002003'
$ws.Range('E51').Value = '-0.03%'

# Reset the style index back to the default ("Normal") now that the
# values are committed as text, so no stray style id is left on the
# individual cells (keeps styles.xml/cell attrs equivalent to source).
$ws.Range("D2:E51").Style = "Normal"
